$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 2 (Leve Item ID 5489)
$ws.Range("H2").Value = 1100
$ws.Range("J2").Value = 1200
$ws.Range("L2").Value = 1200
$ws.Range("N2").Value = -1426

# Row 111 (Leve Item ID 27768)
$ws.Range("H111").Value = 1790.28
$ws.Range("I111").Value = 5864.5
$ws.Range("J111").Value = 1436
$ws.Range("K111").Value = 17593.5
$ws.Range("L111").Value = 4308
$ws.Range("M111").Value = -14526.5
$ws.Range("N111").Value = -10442

# Row 112 (Leve Item ID 27960)
$ws.Range("H112").Value = 1318.4849
$ws.Range("I112").Value = 700
$ws.Range("J112").Value = 1380.3334
$ws.Range("K112").Value = 2100
$ws.Range("L112").Value = 4141.0002
$ws.Range("M112").Value = -992
$ws.Range("N112").Value = -6357.0002

# Row 113 (Leve Item ID 27775)
$ws.Range("H113").Value = 1400
$ws.Range("I113").Value = 1400
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1400
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 1854
$ws.Range("N113").ClearContents()

# Row 116 (Leve Item ID 27778)
$ws.Range("H116").Value = 6190.5
$ws.Range("I116").Value = 6488.125
$ws.Range("J116").Value = 5000
$ws.Range("K116").Value = 6488.125
$ws.Range("L116").Value = 5000
$ws.Range("M116").Value = -3046.125
$ws.Range("N116").Value = -11884

# Row 118 (Leve Item ID 27958)
$ws.Range("H118").Value = 2669.8965
$ws.Range("I118").Value = 930.7692
$ws.Range("J118").Value = 4082.9375
$ws.Range("K118").Value = 2792.3076
$ws.Range("L118").Value = 12248.8125
$ws.Range("M118").Value = -1135.3076
$ws.Range("N118").Value = -15562.8125

# Row 125 (Leve Item ID 36228)
$ws.Range("H125").Value = 4420
$ws.Range("I125").Value = 10000
$ws.Range("J125").Value = 700
$ws.Range("K125").Value = 90000
$ws.Range("L125").Value = 6300
$ws.Range("M125").Value = -87540
$ws.Range("N125").Value = -11220

# Row 138 (Leve Item ID 44169)
$ws.Range("H138").Value = 5210070.5
$ws.Range("I138").Value = 1327.625
$ws.Range("J138").Value = 10418813
$ws.Range("K138").Value = 3982.875
$ws.Range("L138").Value = 31256439
$ws.Range("M138").Value = 1157.125
$ws.Range("N138").Value = -31266719

# Row 139 (Leve Item ID 42306)
$ws.Range("H139").Value = 29316.666
$ws.Range("J139").Value = 29316.666
$ws.Range("L139").Value = 29316.666
$ws.Range("N139").Value = -39596.666

# Row 140 (Leve Item ID 42459)
$ws.Range("H140").Value = 79060
$ws.Range("J140").Value = 79060
$ws.Range("L140").Value = 79060
$ws.Range("N140").Value = -89420

$ws = $wb.Worksheets.Item("ARM")
# Row 32 (Leve Item ID 44147)
$ws.Range("H32").Value = 12815
$ws.Range("I32").Value = 15128.536
$ws.Range("K32").Value = 15128.536
$ws.Range("M32").Value = -14841.536

# Row 61 (Leve Item ID 43999)
$ws.Range("H61").Value = 20837250
$ws.Range("I61").Value = 23813238
$ws.Range("K61").Value = 23813238
$ws.Range("M61").Value = -23813026

# Row 74 (Leve Item ID 44000)
$ws.Range("H74").Value = 15628260
$ws.Range("I74").Value = 27779764
$ws.Range("J74").Value = 4898.2144
$ws.Range("K74").Value = 27779764
$ws.Range("L74").Value = 4898.2144
$ws.Range("M74").Value = -27778890
$ws.Range("N74").Value = -6646.2144

# Row 77 (Leve Item ID 44000)
$ws.Range("H77").Value = 15628260
$ws.Range("I77").Value = 27779764
$ws.Range("J77").Value = 4898.2144
$ws.Range("K77").Value = 138898820
$ws.Range("L77").Value = 24491.072
$ws.Range("M77").Value = -138894452
$ws.Range("N77").Value = -33227.072

# Row 110 (Leve Item ID 27708)
$ws.Range("H110").Value = 1475.05
$ws.Range("I110").Value = 1364.7646
$ws.Range("J110").Value = 2100
$ws.Range("K110").Value = 1364.7646
$ws.Range("L110").Value = 2100
$ws.Range("M110").Value = 680.2354
$ws.Range("N110").Value = -6190

# Row 132 (Leve Item ID 43997)
$ws.Range("H132").Value = 8335940
$ws.Range("I132").Value = 10002173
$ws.Range("K132").Value = 30006519
$ws.Range("M132").Value = -30003989

# Row 136 (Leve Item ID 43999)
$ws.Range("H136").Value = 20837250
$ws.Range("I136").Value = 23813238
$ws.Range("K136").Value = 71439714
$ws.Range("M136").Value = -71437164

$ws = $wb.Worksheets.Item("BSM")
# Row 117 (Leve Item ID 26124)
$ws.Range("H117").Value = 59999.5
$ws.Range("J117").Value = 59999.5
$ws.Range("L117").Value = 59999.5
$ws.Range("N117").Value = -69177.5

# Row 134 (Leve Item ID 43998)
$ws.Range("H134").Value = 4509.067
$ws.Range("I134").Value = 2694.182
$ws.Range("J134").Value = 9500
$ws.Range("K134").Value = 8082.545999999999
$ws.Range("L134").Value = 28500
$ws.Range("M134").Value = -5547.545999999999
$ws.Range("N134").Value = -33570

$ws = $wb.Worksheets.Item("CRP")
# Row 140 (Leve Item ID 42455)
$ws.Range("H140").Value = 35350
$ws.Range("J140").Value = 35350
$ws.Range("L140").Value = 35350
$ws.Range("N140").Value = -45710

$ws = $wb.Worksheets.Item("CUL")
# Row 2 (Leve Item ID 4847)
$ws.Range("H2").Value = 107
$ws.Range("I2").Value = 147
$ws.Range("J2").Value = 27
$ws.Range("K2").Value = 882
$ws.Range("L2").Value = 162
$ws.Range("M2").Value = -769
$ws.Range("N2").Value = -388

# Row 80 (Leve Item ID 12890)
$ws.Range("H80").Value = 3600
$ws.Range("I80").Value = 1033.3334
$ws.Range("J80").Value = 4700
$ws.Range("K80").Value = 3100.0002
$ws.Range("L80").Value = 14100
$ws.Range("M80").Value = -2164.0002
$ws.Range("N80").Value = -15972

# Row 83 (Leve Item ID 12890)
$ws.Range("H83").Value = 3600
$ws.Range("I83").Value = 1033.3334
$ws.Range("J83").Value = 4700
$ws.Range("K83").Value = 9300.000599999999
$ws.Range("L83").Value = 42300
$ws.Range("M83").Value = -4620.000599999999
$ws.Range("N83").Value = -51660

# Row 106 (Leve Item ID 19819)
$ws.Range("H106").Value = 8479.799999999999
$ws.Range("J106").Value = 8479.799999999999
$ws.Range("L106").Value = 25439.4
$ws.Range("N106").Value = -27331.4

# Row 121 (Leve Item ID 27878)
$ws.Range("H121").Value = 1304.3182
$ws.Range("I121").Value = 433.33334
$ws.Range("J121").Value = 1441.8422
$ws.Range("K121").Value = 1300.00002
$ws.Range("L121").Value = 4325.5266
$ws.Range("M121").Value = 9.99998000000005
$ws.Range("N121").Value = -6945.5266

# Row 125 (Leve Item ID 36043)
$ws.Range("H125").Value = 5266.6665
$ws.Range("J125").Value = 5266.6665
$ws.Range("L125").Value = 15799.9995
$ws.Range("N125").Value = -25639.9995

# Row 129 (Leve Item ID 36054)
$ws.Range("H129").Value = 2594.4614
$ws.Range("J129").Value = 3092.4211
$ws.Range("L129").Value = 9277.263300000001
$ws.Range("N129").Value = -19277.2633

$ws = $wb.Worksheets.Item("GSM")
# Row 113 (Leve Item ID 27710)
$ws.Range("H113").Value = 126015.25
$ws.Range("I113").Value = 143946
$ws.Range("J113").Value = 500
$ws.Range("K113").Value = 143946
$ws.Range("L113").Value = 500
$ws.Range("M113").Value = -141776
$ws.Range("N113").Value = -4840

# Row 132 (Leve Item ID 44008)
$ws.Range("H132").Value = 4810.143
$ws.Range("I132").Value = 3803.5557
$ws.Range("J132").Value = 5875.9414
$ws.Range("K132").Value = 11410.6671
$ws.Range("L132").Value = 17627.8242
$ws.Range("M132").Value = -8880.667099999999
$ws.Range("N132").Value = -22687.8242

$ws = $wb.Worksheets.Item("LTW")
# Row 7 (Leve Item ID 36249)
$ws.Range("H7").Value = 6194.6665
$ws.Range("I7").Value = 6510.4
$ws.Range("J7").Value = 5800
$ws.Range("K7").Value = 6510.4
$ws.Range("L7").Value = 5800
$ws.Range("M7").Value = -6398.4
$ws.Range("N7").Value = -6024

# Row 16 (Leve Item ID 5289)
$ws.Range("H16").Value = 2963.2856
$ws.Range("I16").Value = 2711.45
$ws.Range("J16").Value = 8000
$ws.Range("K16").Value = 2711.45
$ws.Range("L16").Value = 8000
$ws.Range("M16").Value = -2541.45
$ws.Range("N16").Value = -8340

# Row 61 (Leve Item ID 27740)
$ws.Range("H61").Value = 1375.25
$ws.Range("I61").Value = 1000
$ws.Range("J61").Value = 1500.3334
$ws.Range("K61").Value = 1000
$ws.Range("L61").Value = 1500.3334
$ws.Range("M61").Value = -798
$ws.Range("N61").Value = -1904.3334

# Row 113 (Leve Item ID 27740)
$ws.Range("H113").Value = 1375.25
$ws.Range("I113").Value = 1000
$ws.Range("J113").Value = 1500.3334
$ws.Range("K113").Value = 1000
$ws.Range("L113").Value = 1500.3334
$ws.Range("M113").Value = 1170
$ws.Range("N113").Value = -5840.3334

# Row 126 (Leve Item ID 36249)
$ws.Range("H126").Value = 6194.6665
$ws.Range("I126").Value = 6510.4
$ws.Range("J126").Value = 5800
$ws.Range("K126").Value = 19531.2
$ws.Range("L126").Value = 17400
$ws.Range("M126").Value = -17061.2
$ws.Range("N126").Value = -22340

# Row 139 (Leve Item ID 43310)
$ws.Range("H139").Value = 39379.4
$ws.Range("J139").Value = 39238.223
$ws.Range("L139").Value = 39238.223
$ws.Range("N139").Value = -49518.223

$ws = $wb.Worksheets.Item("WVR")
# Row 132 (Leve Item ID 44029)
$ws.Range("H132").Value = 2140.625
$ws.Range("I132").Value = 1202
$ws.Range("J132").Value = 2870.6667
$ws.Range("K132").Value = 3606
$ws.Range("L132").Value = 8612.000100000001
$ws.Range("M132").Value = -1076
$ws.Range("N132").Value = -13672.0001

# Row 136 (Leve Item ID 44031)
$ws.Range("H136").Value = 894.3333
$ws.Range("I136").Value = 893.6799999999999
$ws.Range("J136").Value = 902.5
$ws.Range("K136").Value = 2681.04
$ws.Range("L136").Value = 2707.5
$ws.Range("M136").Value = -131.04
$ws.Range("N136").Value = -7807.5

# Row 138 (Leve Item ID 42347)
$ws.Range("H138").Value = 63244.75
$ws.Range("J138").Value = 63244.75
$ws.Range("L138").Value = 63244.75
$ws.Range("N138").Value = -73524.75
